# Actualización automática 2025-06-24 16:45:08
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 (OTROS): VENTA increases, POR CUMPLIR mirrors the negative of VENTA
$ws.Range("D2").Value = 3819.95
$ws.Range("E2").Value = -3819.95

# Row 4 (TOTAL): recompute dependent totals
$ws.Range("D4").Value = 4141.8
$ws.Range("E4").Value = 13358.2
$ws.Range("F4").Value = 0.2366742857142857
